$wb = $excel.ActiveWorkbook

# --- Update the Logs sheet: append new row 36 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(36, 1).Value = "CE-certificaten verzoek"
$logs.Cells.Item(36, 2).Value = "inkoop@testbedrijf123.nl"
$logs.Cells.Item(36, 3).Value = "Kun je mij de CE-certificaten van de EcoPro-700 sturen?"
$logs.Cells.Item(36, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(36, 5).Value = "Bedankt, we hebben dit doorgestuurd naar kwaliteit@testbedrijf123.nl."
$logs.Cells.Item(36, 6).Value = "2025-08-14 22:02:10"
$logs.Cells.Item(36, 7).Value = "Nee"
$logs.Cells.Item(36, 8).Value = "Ja"
$logs.Cells.Item(36, 9).Value = "Nee"
$logs.Cells.Item(36, 10).Value = "Nee"

# --- Extend conditional-formatting ranges to cover the new row ---
$logs.Range("D2:D35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D36"))
$logs.Range("G2:G35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G36"))
$logs.Range("H2:H35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H36"))
$logs.Range("I2:I35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I36"))
$logs.Range("J2:J35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J36"))

# --- Update the Dashboard sheet: bump the count for the "Intern verzoek / Actie voor medewerker" category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 28
